$d = $word.ActiveDocument

$d.Content.Find.Execute("DOCX, DOC, PDF, HTML, XPS, RTF and TXT", $true, $false, $false, $false, $false,
                         $true, 1, $false, "DOCX, DOC, PDF, HTML, XPS, RTF and TXT", 2)
